$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("02-10-2021", "03-10-2021", "04-10-2021", "05-10-2021", "06-10-2021")

$startRow = 276
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
    $cell.ClearFormats()
    $ws.Cells.Item($r, 2).Value = 3623
    $ws.Cells.Item($r, 3).Value = 240
}
